# "Actualizar menu desde Excel" — append a new product ("Tallarin saltado")
# to the Products table, widen the Description/Category columns to fit the
# new (and existing) content, and leave the UI focused back on the Products
# sheet at the freshly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Grow the "products" Excel table by one row (this also advances the table's
# ref/autoFilter range and the sheet dimension automatically).
$lo = $ws.ListObjects.Item("products")
$newRow = $lo.ListRows.Add()

# Fill in the new row (row 33): ID, Name, Description, Price, Category,
# Subcategory, Image, Status.
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Tallarin saltado"
$ws.Range("C33").Value = "Este tallarin es rico"
$ws.Range("D33").Value = 14
$ws.Range("E33").Value = "Almuerzos"
$ws.Range("F33").Value = "Menú del Día"
$ws.Range("G33").Value = "/img/proceso.webp"
$ws.Range("H33").Value = "ACTIVO"

# Extend the Category / Subcategory dropdown validations down to the new row.
$eRange = $ws.Range("E2:E33")
$eRange.Validation.Delete() | Out-Null
$eRange.Validation.Add(3, 1, 1, "Category") | Out-Null

$fRange = $ws.Range("F2:F33")
$fRange.Validation.Delete() | Out-Null
$fRange.Validation.Add(3, 1, 1, "INDIRECT(`$E2)") | Out-Null

# Re-fit the Description and Category columns now that the content changed.
$ws.Columns.Item(3).ColumnWidth = 76
$ws.Columns.Item(5).ColumnWidth = 12

# Leave the Categories sheet's own selection parked on K5 ...
$cat = $wb.Worksheets.Item("Categories")
$cat.Activate() | Out-Null
$cat.Range("K5").Select() | Out-Null

# ... but return focus to Products, on the new row's last cell, as the
# active tab/selection for the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("H33").Select() | Out-Null
